$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (shifts existing data rows down by one)
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the new data point
$ws.Cells.Item(2, 1).Value = 0.4420366287231431
$ws.Cells.Item(2, 2).Value = 4.159920692443848
$ws.Cells.Item(2, 3).Value = -1.852656066417694

# The insert copies formatting down from the header row; strip it so the
# new data row stays unstyled like the rest of the data rows.
$ws.Range("A2:C2").ClearFormats()

# Remove the trailing two rows (the old last row shifted to 23, plus
# what is now the duplicated former-second-to-last row at 22) so the
# dataset stays at 20 data rows total.
$ws.Range("A22:C23").EntireRow.Delete()
